$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row: A12 = "Empty Date" label, B12 = empty cell with the same
# date number format as B2 (copy format only, no value).
$ws.Range("A12").Value = "Empty Date"

$ws.Range("B2").Copy()
$ws.Range("B12").PasteSpecial(-4122)

# Clear the clipboard/marching-ants state left behind by Copy().
$excel.CutCopyMode = 0

# Move the active selection to L9, matching the saved view state.
$null = $ws.Range("L9").Select()
